$d = $word.ActiveDocument

# 1) "Automatic setup" command: prepend "sudo " to the voicecommand invocation.
$d.Content.Find.Execute("voicecommand -s -f .commands.conf", $false, $false, $false, $false, $false, `
    $true, 1, $false, "sudo voicecommand -s -f .commands.conf", 2) | Out-Null

# 2) "Config setup" command: prepend "sudo " to the voicecommand invocation.
$enDash = [char]8211
$d.Content.Find.Execute("voicecommand " + $enDash + "e", $false, $false, $false, $false, $false, `
    $true, 1, $false, "sudo voicecommand " + $enDash + "e", 2) | Out-Null

# 3) Add a new "Remove Autostart at boot" section right after the existing
#    "Autostart at boot" section (i.e. after the "sudo update-rc.d voicecommand
#    defaults" paragraph) and before the "Usage" heading.
$r = $d.Content
$r.Find.Execute("sudo update-rc.d voicecommand defaults") | Out-Null
$r.Collapse(0)

$anchorIndex = $r.Paragraphs.Item(1).Index

$r.InsertParagraphAfter()
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()

$d.Paragraphs.Item($anchorIndex + 2).Range.Text = "Remove Autostart at boot"
$d.Paragraphs.Item($anchorIndex + 2).Style = "Heading 2"
$d.Paragraphs.Item($anchorIndex + 3).Range.Text = "sudo update-rc.d -f voicecommand remove"
$d.Paragraphs.Item($anchorIndex + 4).Range.Text = "sudo rm /etc/init.d/voicecommand"

Write-Output "done"
